# Cleanup excel read code
# Add a "Scores_num" worksheet: a values-only (no formulas) snapshot of the
# "Scores" sheet, inserted right between "Scores" and "ScorecardInfo".

$wb = $excel.ActiveWorkbook

$scores = $wb.Worksheets.Item("Scores")

# Select the whole data range on "Scores" (mirrors selecting it before
# copying) and make sure it is the active sheet while we do so.
[void]$scores.Activate()
[void]$scores.Range("A1:H8").Select()

# Insert the new sheet right after "Scores" (so tab order becomes
# Scores, Scores_num, ScorecardInfo) and name it.
$numSheet = $wb.Worksheets.Add($null, $scores)
$numSheet.Name = "Scores_num"

# Copy the Scores data and paste values only, so formulas like =171/200
# become plain numeric constants on the new sheet.
[void]$scores.Range("A1:H8").Copy()
[void]$numSheet.Range("A1").PasteSpecial(-4163)
$excel.CutCopyMode = $false
